{"js": "// Add a new bulleted list item \"Semantisk s\u00f6kning (\u00f6kad komplexitet)\" right\n// after the existing \"Direktredigering?\" list item, matching its list\n// formatting (same ListParagraph style / same numbered list).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Direktredigering?\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Direktredigering?\"');\n}\n\ntarget.list.load(\"id\");\nawait context.sync();\nconst listId = target.list.id;\n\nconst newPara = target.insertParagraph(\n  \"Semantisk s\u00f6kning (\u00f6kad komplexitet)\",\n  Word.InsertLocation.after\n);\nnewPara.style = target.style;\nnewPara.attachToList(listId, 0);\n\nawait context.sync();\n", "ps1": "# Add a new bulleted list item \"Semantisk s\u00f6kning (\u00f6kad komplexitet)\" right\n# after the existing \"Direktredigering?\" list item, matching its list\n# formatting (same ListParagraph style / same numbered list).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Direktredigering?\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph with text 'Direktredigering?'\"\n}\n\n# Insert a new paragraph mark right after the target paragraph - it inherits\n# the source paragraph's formatting (style + list numbering) automatically.\n$target.Range.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.Text = \"Semantisk s\u00f6kning (\u00f6kad komplexitet)\"\n"}
